$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 79.2
$ws.Range("I42").Value = 79.2
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 237.6
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -7.600000000000023
$ws.Range("N42").ClearContents()

$ws.Range("H94").Value = 2598.8
$ws.Range("I94").Value = 2598.8
$ws.Range("K94").Value = 2598.8
$ws.Range("M94").Value = -2147.8

$ws.Range("H97").Value = 1719.8
$ws.Range("J97").Value = 1719.8
$ws.Range("L97").Value = 5159.4
$ws.Range("N97").Value = -6151.4

$ws.Range("H103").Value = 800.6818
$ws.Range("I103").Value = 511.33334
$ws.Range("K103").Value = 1534.00002
$ws.Range("M103").Value = -948.0000199999999

$ws.Range("H137").Value = 3706744.2
$ws.Range("I137").Value = 5265677
$ws.Range("J137").Value = 4278.5
$ws.Range("K137").Value = 15797031
$ws.Range("L137").Value = 12835.5
$ws.Range("M137").Value = -15794481
$ws.Range("N137").Value = -17935.5

$ws.Range("H138").Value = 2384797.8
$ws.Range("I138").Value = 3673.5625
$ws.Range("J138").Value = 3090316
$ws.Range("K138").Value = 11020.6875
$ws.Range("L138").Value = 9270948
$ws.Range("M138").Value = -5880.6875
$ws.Range("N138").Value = -9281228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1005.99
$ws.Range("I32").Value = 884.8977
$ws.Range("J32").Value = 1894
$ws.Range("K32").Value = 884.8977
$ws.Range("L32").Value = 1894
$ws.Range("M32").Value = -597.8977
$ws.Range("N32").Value = -2468

$ws.Range("H45").Value = 2234.4546
$ws.Range("I45").Value = 2368.4285
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2368.4285
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1991.4285
$ws.Range("N45").Value = -2754

$ws.Range("H61").Value = 250500500
$ws.Range("I61").Value = 333667330
$ws.Range("J61").Value = 1000000
$ws.Range("K61").Value = 333667330
$ws.Range("L61").Value = 1000000
$ws.Range("M61").Value = -333667118
$ws.Range("N61").Value = -1000424

$ws.Range("H74").Value = 10082837
$ws.Range("I74").Value = 13946162
$ws.Range("J74").Value = 148571.42
$ws.Range("K74").Value = 13946162
$ws.Range("L74").Value = 148571.42
$ws.Range("M74").Value = -13945288
$ws.Range("N74").Value = -150319.42

$ws.Range("H77").Value = 10082837
$ws.Range("I77").Value = 13946162
$ws.Range("J77").Value = 148571.42
$ws.Range("K77").Value = 69730810
$ws.Range("L77").Value = 742857.1000000001
$ws.Range("M77").Value = -69726442
$ws.Range("N77").Value = -751593.1000000001

$ws.Range("H102").Value = 8405332
$ws.Range("I102").Value = 11906504
$ws.Range("J102").Value = 2520
$ws.Range("K102").Value = 11906504
$ws.Range("L102").Value = 2520
$ws.Range("M102").Value = -11904882
$ws.Range("N102").Value = -5764

$ws.Range("H110").Value = 1176.4783
$ws.Range("I110").Value = 1189.4286
$ws.Range("J110").Value = 1156.3334
$ws.Range("K110").Value = 1189.4286
$ws.Range("L110").Value = 1156.3334
$ws.Range("M110").Value = 855.5714
$ws.Range("N110").Value = -5246.3334

$ws.Range("H136").Value = 250500500
$ws.Range("I136").Value = 333667330
$ws.Range("J136").Value = 1000000
$ws.Range("K136").Value = 1001001990
$ws.Range("L136").Value = 3000000
$ws.Range("M136").Value = -1000999440
$ws.Range("N136").Value = -3005100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1115.129
$ws.Range("I99").Value = 1160
$ws.Range("J99").Value = 1067.2667
$ws.Range("K99").Value = 1160
$ws.Range("L99").Value = 1067.2667
$ws.Range("M99").Value = 338
$ws.Range("N99").Value = -4063.2667

$ws.Range("H134").Value = 3908.16
$ws.Range("I134").Value = 3616.524
$ws.Range("J134").Value = 5439.25
$ws.Range("K134").Value = 10849.572
$ws.Range("L134").Value = 16317.75
$ws.Range("M134").Value = -8314.572
$ws.Range("N134").Value = -21387.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2065.1
$ws.Range("I16").Value = 1864.4286
$ws.Range("J16").Value = 2533.3333
$ws.Range("K16").Value = 1864.4286
$ws.Range("L16").Value = 2533.3333
$ws.Range("M16").Value = -1577.4286
$ws.Range("N16").Value = -3107.3333

$ws.Range("H31").Value = 2584.449
$ws.Range("I31").Value = 1605.3226
$ws.Range("J31").Value = 4270.722
$ws.Range("K31").Value = 1605.3226
$ws.Range("L31").Value = 4270.722
$ws.Range("M31").Value = -1310.3226
$ws.Range("N31").Value = -4860.722

$ws.Range("H34").Value = 2584.449
$ws.Range("I34").Value = 1605.3226
$ws.Range("J34").Value = 4270.722
$ws.Range("K34").Value = 1605.3226
$ws.Range("L34").Value = 4270.722
$ws.Range("M34").Value = -1403.3226
$ws.Range("N34").Value = -4674.722

$ws.Range("H99").Value = 3892.4348
$ws.Range("I99").Value = 2874.4707
$ws.Range("J99").Value = 6776.6665
$ws.Range("K99").Value = 2874.4707
$ws.Range("L99").Value = 6776.6665
$ws.Range("M99").Value = -1376.4707
$ws.Range("N99").Value = -9772.666499999999

$ws.Range("H107").Value = 607.36365
$ws.Range("I107").Value = 563.875
$ws.Range("J107").Value = 723.3333
$ws.Range("K107").Value = 563.875
$ws.Range("L107").Value = 723.3333
$ws.Range("M107").Value = 1356.125
$ws.Range("N107").Value = -4563.3333

$ws.Range("H113").Value = 2065.1
$ws.Range("I113").Value = 1864.4286
$ws.Range("J113").Value = 2533.3333
$ws.Range("K113").Value = 1864.4286
$ws.Range("L113").Value = 2533.3333
$ws.Range("M113").Value = 305.5714
$ws.Range("N113").Value = -6873.3333

$ws.Range("H126").Value = 3892.4348
$ws.Range("I126").Value = 2874.4707
$ws.Range("J126").Value = 6776.6665
$ws.Range("K126").Value = 8623.4121
$ws.Range("L126").Value = 20329.9995
$ws.Range("M126").Value = -6153.4121
$ws.Range("N126").Value = -25269.9995

$ws.Range("H132").Value = 61730.47
$ws.Range("I132").Value = 2984
$ws.Range("J132").Value = 145654
$ws.Range("K132").Value = 8952
$ws.Range("L132").Value = 436962
$ws.Range("M132").Value = -6422
$ws.Range("N132").Value = -442022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 350.8
$ws.Range("I98").Value = 66.666664
$ws.Range("J98").Value = 777
$ws.Range("K98").Value = 199.999992
$ws.Range("L98").Value = 2331
$ws.Range("M98").Value = 1298.000008
$ws.Range("N98").Value = -5327

$ws.Range("H131").Value = 1186.3889
$ws.Range("I131").Value = 675
$ws.Range("J131").Value = 1227.3
$ws.Range("K131").Value = 2025
$ws.Range("L131").Value = 3681.9
$ws.Range("M131").Value = 3015
$ws.Range("N131").Value = -13761.9

$ws.Range("H132").Value = 1526.6666
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1526.6666
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 13739.9994
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -18799.9994

$ws.Range("H141").Value = 10000
$ws.Range("I141").Value = 10000
$ws.Range("K141").Value = 30000
$ws.Range("M141").Value = -24820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26365.4
$ws.Range("I70").Value = 35956.258
$ws.Range("J70").Value = 5128.5
$ws.Range("K70").Value = 35956.258
$ws.Range("L70").Value = 5128.5
$ws.Range("M70").Value = -35686.258
$ws.Range("N70").Value = -5668.5

$ws.Range("H73").Value = 26365.4
$ws.Range("I73").Value = 35956.258
$ws.Range("J73").Value = 5128.5
$ws.Range("K73").Value = 35956.258
$ws.Range("L73").Value = 5128.5
$ws.Range("M73").Value = -35020.258
$ws.Range("N73").Value = -7000.5

$ws.Range("H80").Value = 3706.1904
$ws.Range("I80").Value = 3465
$ws.Range("J80").Value = 3802.6667
$ws.Range("K80").Value = 3465
$ws.Range("L80").Value = 3802.6667
$ws.Range("M80").Value = -2467
$ws.Range("N80").Value = -5798.6667

$ws.Range("H83").Value = 3706.1904
$ws.Range("I83").Value = 3465
$ws.Range("J83").Value = 3802.6667
$ws.Range("K83").Value = 17325
$ws.Range("L83").Value = 19013.3335
$ws.Range("M83").Value = -12333
$ws.Range("N83").Value = -28997.3335

$ws.Range("H97").Value = 1553.2632
$ws.Range("I97").Value = 1690.8125
$ws.Range("J97").Value = 819.6667
$ws.Range("K97").Value = 1690.8125
$ws.Range("L97").Value = 819.6667
$ws.Range("M97").Value = -1194.8125
$ws.Range("N97").Value = -1811.6667

$ws.Range("H107").Value = 405.3
$ws.Range("I107").Value = 245
$ws.Range("J107").Value = 445.375
$ws.Range("K107").Value = 245
$ws.Range("L107").Value = 445.375
$ws.Range("M107").Value = 1675
$ws.Range("N107").Value = -4285.375

$ws.Range("H122").Value = 2717.6511
$ws.Range("I122").Value = 2430.4194
$ws.Range("J122").Value = 3459.6667
$ws.Range("K122").Value = 7291.2582
$ws.Range("L122").Value = 10379.0001
$ws.Range("M122").Value = -4841.2582
$ws.Range("N122").Value = -15279.0001

$ws.Range("H126").Value = 2435.3333
$ws.Range("I126").Value = 2439.6365
$ws.Range("K126").Value = 7318.9095
$ws.Range("M126").Value = -4848.9095

$ws.Range("H132").Value = 66320.35000000001
$ws.Range("I132").Value = 41717.32
$ws.Range("J132").Value = 168833
$ws.Range("K132").Value = 125151.96
$ws.Range("L132").Value = 506499
$ws.Range("M132").Value = -122621.96
$ws.Range("N132").Value = -511559

$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -95070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2833
$ws.Range("I7").Value = 2600
$ws.Range("J7").Value = 2949.5
$ws.Range("K7").Value = 2600
$ws.Range("L7").Value = 2949.5
$ws.Range("M7").Value = -2488
$ws.Range("N7").Value = -3173.5

$ws.Range("H100").Value = 1323.8108
$ws.Range("I100").Value = 1094.826
$ws.Range("K100").Value = 1094.826
$ws.Range("M100").Value = -553.826

$ws.Range("H126").Value = 2833
$ws.Range("I126").Value = 2600
$ws.Range("J126").Value = 2949.5
$ws.Range("K126").Value = 7800
$ws.Range("L126").Value = 8848.5
$ws.Range("M126").Value = -5330
$ws.Range("N126").Value = -13788.5

$ws.Range("H132").Value = 119267.92
$ws.Range("I132").Value = 74769.28999999999
$ws.Range("J132").Value = 171183
$ws.Range("K132").Value = 224307.87
$ws.Range("L132").Value = 513549
$ws.Range("M132").Value = -221777.87
$ws.Range("N132").Value = -518609

$ws.Range("H136").Value = 83640.44
$ws.Range("I136").Value = 64362.938
$ws.Range("K136").Value = 193088.814
$ws.Range("M136").Value = -190538.814

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 168240.25
$ws.Range("I132").Value = 112542.664
$ws.Range("J132").Value = 335333
$ws.Range("K132").Value = 337627.992
$ws.Range("L132").Value = 1005999
$ws.Range("M132").Value = -335097.992
$ws.Range("N132").Value = -1011059

$ws.Range("H136").Value = 68582.37
$ws.Range("I136").Value = 53656.367
$ws.Range("K136").Value = 160969.101
$ws.Range("M136").Value = -158419.101
